# Applies the "vendas1.xlsx" update: a bunch of corrected K/L (and a couple
# H/I/E/J) values in existing rows 426-542, plus 11 brand-new rows (544-554)
# appended after the previous last row (543). The defined name "Vendas" and
# the sheet dimension both grow from N543 to N554 automatically on save.
#
# NOTE: the row-append logic deliberately avoids a PowerShell function with
# parameters (`function Foo { param(...) ... }`) wrapping Copy()/PasteSpecial()
# calls - combining that pattern with an *earlier* top-level Copy()/PasteSpecial()
# pair triggers a severe slowdown in this COM host. A plain `foreach` loop over
# hashtables has none of that cost and is used instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 426: add a previously-missing "Data_recebida" (column L) value.
# The cell doesn't exist yet, so first clone the date-format style (s="1")
# from a neighbouring L cell, then write the value.
# ---------------------------------------------------------------------------
$ws.Range("L427").Copy()
$ws.Range("L426").PasteSpecial(-4122)
$ws.Cells.Item(426, 12).Value = 44681

# ---------------------------------------------------------------------------
# Simple column L ("Data_recebida") corrections - cells already exist with
# the right style, so a plain value assignment is enough.
# ---------------------------------------------------------------------------
$ws.Cells.Item(505, 12).Value = 44681
$ws.Cells.Item(509, 12).Value = 44681
$ws.Cells.Item(530, 12).Value = 44681
$ws.Cells.Item(536, 12).Value = 44738
$ws.Cells.Item(538, 12).Value = 44681
$ws.Cells.Item(539, 12).Value = 44681
$ws.Cells.Item(542, 12).Value = 44708

# ---------------------------------------------------------------------------
# Row 526: Valor Total / Valor Projeto corrected, and the receipt date moved.
# ---------------------------------------------------------------------------
$ws.Cells.Item(526, 8).Value = 42420.51
$ws.Cells.Item(526, 9).Value = 42420.51
$ws.Cells.Item(526, 12).Value = 44926

# ---------------------------------------------------------------------------
# Row 533: now marked paid ("S"), with the received amount and date filled.
# ---------------------------------------------------------------------------
$ws.Cells.Item(533, 10).Value = "S"
$ws.Cells.Item(533, 11).Value = 157480
$ws.Cells.Item(533, 12).Value = 44677

# ---------------------------------------------------------------------------
# Rows 534 / 535: received amount filled in, receipt date corrected.
# ---------------------------------------------------------------------------
$ws.Cells.Item(534, 11).Value = 64450
$ws.Cells.Item(534, 12).Value = 44670

$ws.Cells.Item(535, 11).Value = 30000
$ws.Cells.Item(535, 12).Value = 44670

# ---------------------------------------------------------------------------
# Row 540: description text expanded, amounts revised upward, now marked
# paid ("S"), and the receipt date corrected.
# ---------------------------------------------------------------------------
$ws.Cells.Item(540, 5).Value = "2 BALCÕES CASTAS E 4 CAIXOTES"
$ws.Cells.Item(540, 8).Value = 3498.81
$ws.Cells.Item(540, 9).Value = 3498.81
$ws.Cells.Item(540, 10).Value = "S"
$ws.Cells.Item(540, 11).Value = 3498.81
$ws.Cells.Item(540, 12).Value = 44650

# ---------------------------------------------------------------------------
# Row 541: amounts revised upward, receipt date corrected.
# ---------------------------------------------------------------------------
$ws.Cells.Item(541, 8).Value = 4500
$ws.Cells.Item(541, 9).Value = 4500
$ws.Cells.Item(541, 11).Value = 4500
$ws.Cells.Item(541, 12).Value = 44655

# ---------------------------------------------------------------------------
# Append 11 brand-new sales rows (544-554). Clone formatting from the row
# immediately above with Copy/PasteSpecial so every column keeps the exact
# same style indices already used throughout the sheet (s="1" dates, s="2"
# numbers, default for text/ints), then fill in the values.
# ---------------------------------------------------------------------------
$newRows = @(
    @{ Row=544; A=543; B=44644; D="V3A"; E="PETROBRAS ROVR NA RIO2C"; H=269000; I=269000; J="N"; K=0; L=44718; M=0; N=0 },
    @{ Row=545; A=544; B=44655; D="ÓTIMA CONCESSIONÁRIA"; E="WARNER ANIMAIS FANTÁSTICOS"; H=83000; I=83000; J="N"; K=83000; L=44724; M=0; N=0 },
    @{ Row=546; A=545; B=44655; D="NETZA"; E="STAND EVE RIO2C"; H=64459.99; I=64459.99; J="N"; K=64459.99; L=44718; M=0; N=0 },
    @{ Row=547; A=546; B=44658; D="CERVEJARIA  PRAYA"; E="BAILE DO ENCANTO"; H=50000; I=50000; J="N"; K=50000; L=44681; M=0; N=0 },
    @{ Row=548; A=547; B=44658; D="AGÊNCIA TERRUÁ"; E="STAND BANCO DO BRASIL RIO2C"; H=200000; I=200000; J="N"; K=200000; L=44719; M=0; N=0 },
    @{ Row=549; A=548; B=44664; D="LVHM"; E="REMONTAGEM BAR BELVERDE"; H=7500; I=7500; J="N"; K=7500; L=44681; M=0; N=0 },
    @{ Row=550; A=549; B=44665; D="GLOBO SAT"; E="TÚNEL DO AMOR"; H=300000; I=300000; J="N"; K=300000; L=44712; M=0; N=0 },
    @{ Row=551; A=550; B=44666; D="GLOBO COMUNICAÇÕES"; E="CUBO RIO2C"; H=75000; I=75000; J="N"; K=75000; L=44712; M=0; N=0 },
    @{ Row=552; A=551; B=44667; D="BE COMUNICA"; E="STAND DASA"; H=85248; I=85248; J="N"; K=85248; L=44717; M=0; N=0 },
    @{ Row=553; A=552; B=44683; D="DIALOGO URBANO - BARBARA SOLEDADE"; E="STAND SUBMARINO (OBVIOUS)"; H=20000; I=20000; J="N"; K=20000; L=44742; M=0; N=0 },
    @{ Row=554; A=553; B=44685; D="VOID"; E="VOID TIJUCA"; H=49500; I=49500; J="N"; K=49500; L=44713; M=0; N=0 }
)

foreach ($r in $newRows) {
    $prevRow = $r.Row - 1
    $ws.Range("A" + $prevRow + ":N" + $prevRow).Copy()
    $ws.Range("A" + $r.Row + ":N" + $r.Row).PasteSpecial(-4122)

    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 8).Value = $r.H
    $ws.Cells.Item($r.Row, 9).Value = $r.I
    $ws.Cells.Item($r.Row, 10).Value = $r.J
    $ws.Cells.Item($r.Row, 11).Value = $r.K
    $ws.Cells.Item($r.Row, 12).Value = $r.L
    $ws.Cells.Item($r.Row, 13).Value = $r.M
    $ws.Cells.Item($r.Row, 14).Value = $r.N
}

# ---------------------------------------------------------------------------
# Update the workbook-level defined name "Vendas" to cover the new extent.
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Vendas") {
        $n.RefersTo = "='Vendas'!`$A`$1:`$N`$554"
    }
}
